# The header row (row 1) on Sheet1 gains a new column "Diad1_prom/std_betweendiads"
# inserted before the existing "Diad2_height" column (column O), shifting every
# header from O1..AA1 one column to the right (ending at AB1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank column at O - this pushes O1:AA1 ("Diad2_height" ... "Mean_HB_prom")
# one column to the right, so they now occupy P1:AB1.
$ws.Range("O1").EntireColumn.Insert()

# Match the header formatting used by the rest of row 1 (bold font, thin border,
# centered/top aligned) by copying it from the neighboring header cell.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)

# Set the text for the newly inserted header cell.
$ws.Range("O1").Value = "Diad1_prom/std_betweendiads"
